$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 with the new "custom accuracy" values (rounded to ~2 decimal
# places, per the commit message), including the refined timestamp in A5.
$ws.Range("A5").Value = 44781.9027662037
$ws.Range("B5").Value = 8.220000000000001
$ws.Range("C5").Value = 5.88
$ws.Range("D5").Value = 0.14
$ws.Range("E5").Value = 16.06
$ws.Range("F5").Value = 13.64
$ws.Range("G5").Value = 6.33
$ws.Range("H5").Value = 24.42
$ws.Range("I5").Value = 8.92
$ws.Range("J5").Value = 4.17
$ws.Range("K5").Value = 6.64
$ws.Range("L5").Value = 6.56
$ws.Range("M5").Value = 6.71
$ws.Range("N5").Value = 1.94
$ws.Range("O5").Value = 5.78
$ws.Range("P5").Value = 8.869999999999999
$ws.Range("Q5").Value = 4.84
$ws.Range("R5").Value = 0.3
$ws.Range("S5").Value = 0.14
$ws.Range("T5").Value = 84.05
$ws.Range("U5").Value = 16.85
$ws.Range("V5").Value = 5.72
$ws.Range("W5").Value = 11.47
$ws.Range("X5").Value = 5.81
$ws.Range("Y5").Value = 0.79
$ws.Range("Z5").Value = 11.84
$ws.Range("AA5").Value = 4.88
$ws.Range("AB5").Value = 4.13
$ws.Range("AC5").Value = 4.91
$ws.Range("AD5").Value = 7.37
$ws.Range("AE5").Value = 0.52
$ws.Range("AF5").Value = 21.94
$ws.Range("AG5").Value = 3.01
$ws.Range("AH5").Value = 6.77

# The old row 6 (second 10-minute reading) was dropped entirely, so remove
# it and let Excel shrink the sheet's used range (dimension) accordingly.
$ws.Rows.Item(6).Delete()
